$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header ("Other Informations") - must be entered first so it
# becomes shared-string index 17, matching the target workbook.
$ws.Range("F1").Value = "Other Informations"

# New row 5 - written in the same left-to-right/top-to-bottom order the
# author used so new shared strings land at the expected indices
# (18 "2", 19 "Encounter Table", 20 description, 21 "0212CBB4").
$ws.Range("B5").Value = "2"
$ws.Range("D5").Value = "Encounter Table"
$ws.Range("E5").Value = "The encounter table that rules which monsters can be encountered"
$ws.Range("A5").Value = "0212CBB4"
$ws.Range("C5").Value = "Unsigned"

# Match column F's width to the existing "Description" column (E) - reuse
# the value Excel itself reports for column E's ColumnWidth.
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Restore the selection left behind after the edit.
[void]$ws.Range("A7").Select()
